$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mapping of row -> (DAMSLTag, DialogAct)
$updates = @{
    9  = @("sv", "Statement-opinion")
    23 = @("sv", "Statement-opinion")
    29 = @("sd", "Statement-non-opinion")
    44 = @("aa", "Agree/Accept")
    57 = @("ba", "Appreciation")
    63 = @("sv", "Statement-opinion")
    68 = @("ba", "Appreciation")
    70 = @("sv", "Statement-opinion")
    71 = @("ba", "Appreciation")
    78 = @("sv", "Statement-opinion")
    83 = @("sv", "Statement-opinion")
    87 = @("sv", "Statement-opinion")
    88 = @("aa", "Agree/Accept")
}

foreach ($row in $updates.Keys) {
    $values = $updates[$row]
    $ws.Range("I$row").Value = $values[0]
    $ws.Range("J$row").Value = $values[1]
}
